$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.636.83"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.643.62"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.505"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").Value = "1.872.54"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.02%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.633.72"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.87%  "
$ws.Range("D17").Value = "26.676.42"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("E22").Value = "  +2.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.19%  "
$ws.Range("E24").Value = "  +10.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0517"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.53%  "
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("D34").Value = "1.271.97"
$ws.Range("E34").Value = "  +4.78%  "
$ws.Range("E35").Value = "  +2.39%  "
$ws.Range("E36").Value = "  +5.58%  "
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.531"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.829"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.813"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.62%  "
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.28%  "
$ws.Range("D44").Value = "1.782.59"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.61%  "
$ws.Range("E47").Value = "  +3.19%  "
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0975"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.408"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.31%  "
